# Updates the cryptos list (Price column D and Volume(1h) column E) to reflect
# the latest scraped values, including two coin-rank swaps (rows 33/34 and 42/43)
# where the underlying coin data (Coin name + Link) moved to the other row.
#
# Column D ("Price") cells are forced to Text format before assignment so that
# values such as "1.00", "0.0675" or thousands-grouped values like "68.087.78"
# are stored exactly as text (matching the source data) instead of being
# auto-converted to numbers by Excel (which would drop trailing zeros or
# introduce floating point artifacts). The cell style is reset back to
# "Normal" afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "68.087.78"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.48%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.308.73"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.16%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "186.60"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.19%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "583.14"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.25%  "
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 5).Value = "  -1.08%  "
$ws.Cells.Item(9, 5).Value = "  -0.11%  "
$ws.Cells.Item(10, 5).Value = "  +1.44%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.409"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.14%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "3.884.83"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.35%  "
$ws.Cells.Item(13, 5).Value = "  -2.25%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "27.52"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.67%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "68.270.80"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.70%  "
$ws.Cells.Item(16, 5).Value = "  -0.19%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "3.315.59"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "447.34"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +11.47%  "
$ws.Cells.Item(19, 5).Value = "  -0.07%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "13.55"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.70%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "7.79"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.83%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "75.06"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +5.78%  "
$ws.Cells.Item(23, 5).Value = "  +0.02%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "3.459.35"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.49%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.20%  "
$ws.Cells.Item(26, 5).Value = "  +1.06%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.189"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.77%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -4.05%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.44%  "
$ws.Cells.Item(30, 5).Value = "  +1.56%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "22.88"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.05%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "5.36"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.10%  "
$ws.Cells.Item(33, 2).Value = "Fetch.AI"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.29%  "
$ws.Cells.Item(34, 2).Value = "USDe"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.03%  "
$ws.Cells.Item(35, 5).Value = "  -2.00%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.53"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.29%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "163.73"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.04%  "
$ws.Cells.Item(38, 5).Value = "  -2.04%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "27.04"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.24%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "4.51"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.783"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.79%  "
$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.722.38"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.69%  "
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "6.39"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.50%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "40.69"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.12%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.0675"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.53%  "
$ws.Cells.Item(46, 5).Value = "  -0.98%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "24.71"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.58%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "328.16"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.69%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.0277"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.80%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "31.74"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +3.62%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.59%  "
